# Data/UnAuth_Customers.xlsx — refresh the Customer_ID values on Sheet1
# (newer export batch: 177052xx -> 177075xx) and restore the active
# selection on that sheet to A2:C11.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Keep these as text (they were stored as shared-string text originally,
# not numbers), so format the column as Text before writing the new IDs.
$idRange = $ws1.Range("B2:B6")
$idRange.NumberFormat = "@"

$ws1.Cells.Item(2, 2).Value = "17707515"
$ws1.Cells.Item(3, 2).Value = "17707516"
$ws1.Cells.Item(4, 2).Value = "17707517"
$ws1.Cells.Item(5, 2).Value = "17707519"
$ws1.Cells.Item(6, 2).Value = "17707520"

# Update the sheet's active selection to A2:C11.
$ws1.Activate()
$ws1.Range("A2:C11").Select()
